$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OppDealTeamMembers")
$ws.Range("A21").Value = "Laura Kimmel"
$ws.Range("A26").Value = "JP Hanson"
Write-Host "done"
